$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "303.71"
Set-TextValue $ws.Range("E2") "-1.04%"

# Row 3
Set-TextValue $ws.Range("D3") "35.72"
Set-TextValue $ws.Range("E3") "-0.44%"

# Row 4
Set-TextValue $ws.Range("D4") "5.025"
Set-TextValue $ws.Range("E4") "-0.74%"

# Row 5
Set-TextValue $ws.Range("D5") "0.07983"
Set-TextValue $ws.Range("E5") "-1.50%"

# Row 6
Set-TextValue $ws.Range("D6") "1.853"
Set-TextValue $ws.Range("E6") "-4.74%"

# Row 7
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D7") "7.758"
Set-TextValue $ws.Range("E7") "-0.64%"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9211"
Set-TextValue $ws.Range("E8") "-1.50%"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D9") "0.1275"
Set-TextValue $ws.Range("E9") "-3.99%"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1889"
Set-TextValue $ws.Range("E10") "-1.26%"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.09007"
Set-TextValue $ws.Range("E11") "-2.18%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.03418"
Set-TextValue $ws.Range("E12") "-2.90%"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D13") "0.09856"
Set-TextValue $ws.Range("E13") "-0.25%"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D14") "0.001415"
Set-TextValue $ws.Range("E14") "-0.07%"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D15") "0.006311"
Set-TextValue $ws.Range("E15") "9.02%"

# Row 16
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D16") "3.858"
Set-TextValue $ws.Range("E16") "7.12%"

# Row 17
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D17") "4.118"
Set-TextValue $ws.Range("E17") "-0.64%"

# Row 18
Set-TextValue $ws.Range("E18") "14.23%"

# Row 19
Set-TextValue $ws.Range("D19") "0.3404"
Set-TextValue $ws.Range("E19") "-0.73%"

# Row 20
Set-TextValue $ws.Range("E20") "0.49%"

# Row 21
Set-TextValue $ws.Range("D21") "4.798"
Set-TextValue $ws.Range("E21") "-7.47%"

# Row 22
Set-TextValue $ws.Range("D22") "0.2339"
Set-TextValue $ws.Range("E22") "-10.60%"

# Row 23
Set-TextValue $ws.Range("E23") "-0.87%"

# Row 24
Set-TextValue $ws.Range("D24") "0.001234"
Set-TextValue $ws.Range("E24") "0.67%"

# Row 25
Set-TextValue $ws.Range("D25") "0.004848"
Set-TextValue $ws.Range("E25") "1.55%"

# Row 27
Set-TextValue $ws.Range("E27") "-0.33%"

# Row 28
Set-TextValue $ws.Range("E28") "42.18%"

# Row 39
Set-TextValue $ws.Range("D39") "0.01926"
Set-TextValue $ws.Range("E39") "-3.35%"

# Row 40
Set-TextValue $ws.Range("D40") "0.05117"
Set-TextValue $ws.Range("E40") "2.16%"

# Row 41
Set-TextValue $ws.Range("D41") "0.007561"
Set-TextValue $ws.Range("E41") "-0.89%"

# Row 42
Set-TextValue $ws.Range("D42") "0.01011"
Set-TextValue $ws.Range("E42") "-10.01%"

# Row 43
Set-TextValue $ws.Range("D43") "0.1344"
Set-TextValue $ws.Range("E43") "-2.50%"

# Row 44
Set-TextValue $ws.Range("E44") "0.14%"

# Row 45
Set-TextValue $ws.Range("D45") "0.009850"
Set-TextValue $ws.Range("E45") "-13.42%"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006175"
Set-TextValue $ws.Range("E46") "-3.42%"

# Row 47
Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "-0.38%"

# Row 48
Set-TextValue $ws.Range("D48") "63.68"
Set-TextValue $ws.Range("E48") "0.19%"

# Row 49
Set-TextValue $ws.Range("D49") "0.001250"
Set-TextValue $ws.Range("E49") "5.03%"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002100"
Set-TextValue $ws.Range("E50") "-0.38%"

# Row 51
Set-TextValue $ws.Range("D51") "0.0002000"
Set-TextValue $ws.Range("E51") "-0.38%"
